$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1588.8
$ws.Range("I34").Value = 1588.8
$ws.Range("K34").Value = 1588.8
$ws.Range("M34").Value = -1385.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 1588.8
$ws.Range("I36").Value = 1588.8
$ws.Range("K36").Value = 1588.8
$ws.Range("M36").Value = -873.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 11497050
$ws.Range("I76").Value = 2854.1177
$ws.Range("J76").Value = 27780494
$ws.Range("K76").Value = 2854.1177
$ws.Range("L76").Value = 27780494
$ws.Range("M76").Value = -2539.1177
$ws.Range("N76").Value = -27781124

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 11497050
$ws.Range("I79").Value = 2854.1177
$ws.Range("J79").Value = 27780494
$ws.Range("K79").Value = 2854.1177
$ws.Range("L79").Value = 27780494
$ws.Range("M79").Value = -1762.1177
$ws.Range("N79").Value = -27782678

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 303.34616
$ws.Range("I92").Value = 289.38095
$ws.Range("K92").Value = 289.38095
$ws.Range("M92").Value = 958.61905

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 887.25714
$ws.Range("I137").Value = 745.8475
$ws.Range("J137").Value = 1645.7273
$ws.Range("K137").Value = 2237.5425
$ws.Range("L137").Value = 4937.1819
$ws.Range("M137").Value = 312.4575
$ws.Range("N137").Value = -10037.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14929463
$ws.Range("I32").Value = 4305.1816
$ws.Range("J32").Value = 83336430
$ws.Range("K32").Value = 4305.1816
$ws.Range("L32").Value = 83336430
$ws.Range("M32").Value = -4018.1816
$ws.Range("N32").Value = -83337004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1761.0769
$ws.Range("I63").Value = 1761.0769
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1761.0769
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1075.0769
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1761.0769
$ws.Range("I66").Value = 1761.0769
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8805.3845
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5373.3845
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1356.3334
$ws.Range("I122").Value = 1253.5714
$ws.Range("K122").Value = 3760.7142
$ws.Range("M122").Value = -1310.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 8496
$ws.Range("I75").Value = 8496
$ws.Range("K75").Value = 8496
$ws.Range("M75").Value = -7560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H78").Value = 8496
$ws.Range("I78").Value = 8496
$ws.Range("K78").Value = 25488
$ws.Range("M78").Value = -20808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 17210.334
$ws.Range("I33").Value = 1631
$ws.Range("K33").Value = 1631
$ws.Range("M33").Value = -1252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 30128
$ws.Range("J80").Value = 30128
$ws.Range("L80").Value = 30128
$ws.Range("N80").Value = -32374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 70000
$ws.Range("J81").Value = 70000
$ws.Range("L81").Value = 70000
$ws.Range("N81").Value = -71996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 18266.334
$ws.Range("J82").Value = 18266.334
$ws.Range("L82").Value = 18266.334
$ws.Range("N82").Value = -18988.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 30128
$ws.Range("J83").Value = 30128
$ws.Range("L83").Value = 90384
$ws.Range("N83").Value = -101616

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 70000
$ws.Range("J84").Value = 70000
$ws.Range("L84").Value = 210000
$ws.Range("N84").Value = -219984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 18266.334
$ws.Range("J85").Value = 18266.334
$ws.Range("L85").Value = 18266.334
$ws.Range("N85").Value = -20762.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 27335.75
$ws.Range("J88").Value = 27335.75
$ws.Range("L88").Value = 27335.75
$ws.Range("N88").Value = -28147.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 27335.75
$ws.Range("J91").Value = 27335.75
$ws.Range("L91").Value = 27335.75
$ws.Range("N91").Value = -30143.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27032.162
$ws.Range("I12").Value = 4.4444447
$ws.Range("J12").Value = 35719.645
$ws.Range("K12").Value = 13.3333341
$ws.Range("L12").Value = 107158.935
$ws.Range("M12").Value = 159.6666659
$ws.Range("N12").Value = -107504.935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 47843.43
$ws.Range("J58").Value = 55666.5
$ws.Range("L58").Value = 166999.5
$ws.Range("N58").Value = -167255.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7100.933
$ws.Range("I68").Value = 401.4
$ws.Range("J68").Value = 20500
$ws.Range("K68").Value = 1204.2
$ws.Range("L68").Value = 61500
$ws.Range("M68").Value = -393.1999999999998
$ws.Range("N68").Value = -63122

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 10634.25
$ws.Range("I69").Value = 350
$ws.Range("J69").Value = 12691.1
$ws.Range("K69").Value = 1050
$ws.Range("L69").Value = 38073.3
$ws.Range("M69").Value = -239
$ws.Range("N69").Value = -39695.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1575
$ws.Range("I70").Value = 906.6667
$ws.Range("J70").Value = 1976
$ws.Range("K70").Value = 2720.0001
$ws.Range("L70").Value = 5928
$ws.Range("M70").Value = -2405.0001
$ws.Range("N70").Value = -6558

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 7100.933
$ws.Range("I71").Value = 401.4
$ws.Range("J71").Value = 20500
$ws.Range("K71").Value = 3612.6
$ws.Range("L71").Value = 184500
$ws.Range("M71").Value = 443.4000000000001
$ws.Range("N71").Value = -192612

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 10634.25
$ws.Range("I72").Value = 350
$ws.Range("J72").Value = 12691.1
$ws.Range("K72").Value = 3150
$ws.Range("L72").Value = 114219.9
$ws.Range("M72").Value = 906
$ws.Range("N72").Value = -122331.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1575
$ws.Range("I73").Value = 906.6667
$ws.Range("J73").Value = 1976
$ws.Range("K73").Value = 2720.0001
$ws.Range("L73").Value = 5928
$ws.Range("M73").Value = -1628.0001
$ws.Range("N73").Value = -8112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 566.3333
$ws.Range("I98").Value = 546.5
$ws.Range("J98").Value = 572
$ws.Range("K98").Value = 1639.5
$ws.Range("L98").Value = 1716
$ws.Range("M98").Value = -141.5
$ws.Range("N98").Value = -4712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 884.63
$ws.Range("J131").Value = 897.567
$ws.Range("L131").Value = 2692.701
$ws.Range("N131").Value = -12772.701

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2811.8845
$ws.Range("I136").Value = 861.2
$ws.Range("J136").Value = 5471.909
$ws.Range("K136").Value = 2583.6
$ws.Range("L136").Value = 16415.727
$ws.Range("M136").Value = 2516.4
$ws.Range("N136").Value = -26615.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 222345.33
$ws.Range("I139").Value = 664.1177
$ws.Range("K139").Value = 1992.3531
$ws.Range("M139").Value = 3147.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6253078
$ws.Range("I80").Value = 3267.8572
$ws.Range("J80").Value = 50001750
$ws.Range("K80").Value = 3267.8572
$ws.Range("L80").Value = 50001750
$ws.Range("M80").Value = -2269.8572
$ws.Range("N80").Value = -50003746

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6253078
$ws.Range("I83").Value = 3267.8572
$ws.Range("J83").Value = 50001750
$ws.Range("K83").Value = 16339.286
$ws.Range("L83").Value = 250008750
$ws.Range("M83").Value = -11347.286
$ws.Range("N83").Value = -250018734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1524.7222
$ws.Range("I102").Value = 1499.6666
$ws.Range("J102").Value = 1650
$ws.Range("K102").Value = 1499.6666
$ws.Range("L102").Value = 1650
$ws.Range("M102").Value = 122.3334
$ws.Range("N102").Value = -4894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4314.718
$ws.Range("I132").Value = 1445.8667
$ws.Range("J132").Value = 13877.556
$ws.Range("K132").Value = 4337.6001
$ws.Range("L132").Value = 41632.66800000001
$ws.Range("M132").Value = -1807.6001
$ws.Range("N132").Value = -46692.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 4933.3335
$ws.Range("J21").Value = 4933.3335
$ws.Range("L21").Value = 4933.3335
$ws.Range("N21").Value = -5281.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 3520
$ws.Range("I107").Value = 3520
$ws.Range("K107").Value = 3520
$ws.Range("M107").Value = -1600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8649.102999999999
$ws.Range("I132").Value = 2742.7646
$ws.Range("K132").Value = 8228.293799999999
$ws.Range("M132").Value = -5698.293799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17751.486
$ws.Range("I132").Value = 20537.611
$ws.Range("J132").Value = 8348.3125
$ws.Range("K132").Value = 61612.833
$ws.Range("L132").Value = 25044.9375
$ws.Range("M132").Value = -59082.833
$ws.Range("N132").Value = -30104.9375
